{"js": "// Merge the three CI/CD description paragraphs into a single updated\n// paragraph describing the SonarCloud-based QA workflow.\nconst body = context.document.body;\n\n// Locate the first paragraph of the three (\"In every pull request, ...\").\nconst results = body.search(\n  \"In every pull request, a github workflow (defined with a yml file in the repo) launches the containers (docker compose) and run the tests.\",\n  { matchCase: false, matchWholeWord: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph not found\");\n}\n\nconst firstParagraph = results.items[0].paragraphs.getFirst();\n\nconst newText =\n  \"In every pull request, a github workflow (defined with a yml file in the repo) runs the tests. We also use SonarCloud for CI: in every pull request, SonarCloud\\u2019s workflow will perform a static code analysis. All the tests and the quality gate must be passed in order for the PR to be accepted (may be accepted/refused automatically or manually).\";\n\n// Replace the whole first paragraph's text (run content) with the new text.\nfirstParagraph.getRange().insertText(newText, \"Replace\");\n\n// The two following paragraphs (\"There is also automatic ...\" and\n// \"The PR may then be accepted ...\") are removed entirely, merging their\n// content into the first paragraph above.\nconst secondParagraph = firstParagraph.getNext();\nconst thirdParagraph = secondParagraph.getNext();\nsecondParagraph.delete();\nthirdParagraph.delete();\n\nawait context.sync();\n", "ps1": "# Merge the three CI/CD description paragraphs into a single updated\n# paragraph describing the SonarCloud-based QA workflow.\n$d = $word.ActiveDocument\n\n# Locate the first of the three paragraphs via Find, so we don't depend on a\n# hard-coded paragraph index.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"In every pull request, a github workflow*run the tests.\", $false, $false, $true)\n\nif (-not $found) {\n    throw \"Target paragraph not found\"\n}\n\n$targetParagraph = $rng.Paragraphs(1)\n$targetRange = $targetParagraph.Range\n\n$newText = \"In every pull request, a github workflow (defined with a yml file in the repo) runs the tests. We also use SonarCloud for CI: in every pull request, SonarCloud\" + [char]0x2019 + \"s workflow will perform a static code analysis. All the tests and the quality gate must be passed in order for the PR to be accepted (may be accepted/refused automatically or manually).\"\n\n# Replace the first paragraph's text, keeping its paragraph mark / formatting.\n$targetRange.Text = $newText\n\n# Remove the two following paragraphs (\"There is also automatic ...\" and\n# \"The PR may then be accepted ...\"), merging their content into the\n# paragraph above.\n$nextParagraph = $targetParagraph.Next()\n$nextParagraph.Range.Delete()\n$nextParagraph2 = $targetParagraph.Next()\n$nextParagraph2.Range.Delete()\n"}
